# Apply the edits described in the commit "minor tweaks, add some maxes"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Opportunity")

# --- Row 22 "Basic Fleet": add Store1Resource "Max 3" and Required2Resource "Max 2" ---
$ws.Range("E22").Value = "Max 3"
$ws.Range("L22").Value = "Max 2"

# --- Row 27 "Expedited Opportunity": remove Consume1Resource (energy) cell entirely ---
$ws.Range("D27").Clear()

# --- Row 28 "New Strategic Directions": Required1Resource officer -> officers ---
$ws.Range("H28").Value = "officers"

# --- Row 32 "Dangerous Mining": add Store1Resource / Required2Resource "Max 5" ---
$ws.Range("E32").Value = "Max 5"
$ws.Range("L32").Value = "Max 5"

# --- Row 34 "Converted Shipyard": add Store1Resource / Required2Resource "Max 3" ---
$ws.Range("E34").Value = "Max 3"
$ws.Range("L34").Value = "Max 3"

# --- Row 40 "Naval Shipyard": add Store1Resource "Max 4" / Required2Resource "Max 3" ---
$ws.Range("E40").Value = "Max 4"
$ws.Range("L40").Value = "Max 3"

# --- Update the active selection on the sheet (was C42, now I27) ---
$ws.Activate() | Out-Null
$ws.Range("I27").Select() | Out-Null
